$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.952.03"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.546.07"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'305.87"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'37.39"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").Value = "'0.0821"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "2.935.17"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "2.531.08"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'15.28"
$ws.Range("E16").Value = "  +7.46%  "
$ws.Range("D17").Value = "'0.878"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "42.974.71"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "'6.56"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'71.76"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'253.90"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").Value = "'27.62"
$ws.Range("E26").Value = "  -4.68%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +9.95%  "
$ws.Range("D29").Value = "'10.24"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'39.00"
$ws.Range("E30").Value = "  +5.12%  "
$ws.Range("D31").Value = "'6.23"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "'157.84"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'0.0801"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").Value = "'3.29"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.65"
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'18.66"
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'24.22"
$ws.Range("E40").Value = "  +5.67%  "
$ws.Range("D41").Value = "'3.48"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'3.90"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'0.0305"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "2.067.87"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'86.33"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "'9.03"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "2.791.32"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").Value = "'103.75"
$ws.Range("E51").Value = "  -1.75%  "
